# "new word and padez"
# Appends 15 new vocabulary rows (332-346) to the word_db sheet:
#   - 4 new "clothes" verbs (rows 332-335)
#   - 11 new "family" nouns (rows 336-346)
# and moves the view/selection down to the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=original, B=translate, D=general_thema_en, E=general_thema_ru,
#          F=main_thema_en, G=main_thema_ru, J=part_of_speech, K=stage,
#          L=lesson, M=type
$newRows = @(
    @{ A="입다";          B="надевать (брюки, юбку...)"; D="action"; E="действие"; F="clothes"; G="одежда";   J="verb"; K=2; L=7; M="word" },
    @{ A="들다";          B="держать (сумку)";           D="action"; E="действие"; F="items";   G="предметы"; J="verb"; K=2; L=7; M="word" },
    @{ A="넥타이를 매다"; B="носить галстук";             D="action"; E="действие"; F="clothes"; G="одежда";   J="verb"; K=2; L=7; M="word" },
    @{ A="목도리를 하다"; B="завязывать шарф";            D="action"; E="действие"; F="clothes"; G="одежда";   J="verb"; K=2; L=7; M="word" },

    @{ A="가족";   B="семья";                        D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="할아버지"; B="дедушка";                      D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="언니";   B="сестра для девушка";             D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="누나";   B="сестро для мальчика";            D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="할머니"; B="бабушка";                        D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="오빠";   B="старший брат для девушка";       D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="형";     B="старший брат для мальчика";      D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="아버지"; B="дедушка";                        D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="동생";   B="младший брат или сетсра";        D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="여동생"; B="младший брат";                   D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" },
    @{ A="남동생"; B="младшая сестра";                 D="people"; E="люди"; F="family"; G="семья"; J="noun"; K=2; L=8; M="word" }
)

$startRow = 332
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value  = $row.A   # A
    $ws.Cells.Item($r, 2).Value  = $row.B   # B
    $ws.Cells.Item($r, 4).Value  = $row.D   # D
    $ws.Cells.Item($r, 5).Value  = $row.E   # E
    $ws.Cells.Item($r, 6).Value  = $row.F   # F
    $ws.Cells.Item($r, 7).Value  = $row.G   # G
    $ws.Cells.Item($r, 10).Value = $row.J   # J
    $ws.Cells.Item($r, 11).Value = $row.K   # K
    $ws.Cells.Item($r, 12).Value = $row.L   # L
    $ws.Cells.Item($r, 13).Value = $row.M   # M
}

# Move the view to the newly-entered rows, as in the saved workbook.
$ws.Range("B338").Select()
$excel.ActiveWindow.ScrollRow = 334
$excel.ActiveWindow.ScrollColumn = 1
